$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Tasks")

# Row 24 (task 8.1): add Status "in progress" in column D
$ws.Range("D24").Value = "in progress"
$ws.Rows.Item(24).RowHeight = 14.9

# Row 27 (task 9.1): add Status "in progress" in column D
$ws.Range("D27").Value = "in progress"

# Move the active selection to D25
$ws.Range("D25").Select()
